# Generate Report for Handback
# The f6dc4ae2-9dd8-455d-8905-b66f00492b6a.md file has now been handed back
# (in sync with en-US) for both the zh-cn and de-de locales. Update the
# Overview sheet and the per-locale status sheets to reflect the new
# status + handback timestamps, and clear the stale "handback file not
# latest" error now that the file is current.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K3").Value = "2016-08-27 04:46:59"
$wsZhCn.Range("P3").Value = ""

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K3").Value = "2016-08-27 04:47:10"
$wsDeDe.Range("P3").Value = ""
